$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowRange($rowA, $rowB, $colStart, $colEnd) {
    for ($col = $colStart; $col -le $colEnd; $col++) {
        $cellA = $ws.Cells.Item($rowA, $col)
        $cellB = $ws.Cells.Item($rowB, $col)
        $valA = $cellA.Value()
        $valB = $cellB.Value()
        $cellA.Value = $valB
        $cellB.Value = $valA
    }
}

# Swap F:V between rows 90 and 91 (match order got corrected)
Swap-RowRange 90 91 6 22

# Swap F:V between rows 122 and 123 (match order got corrected)
Swap-RowRange 122 123 6 22

# Append new rows 125-143 (A=124..142), copying style from row 124
$ws.Range("A124:V124").Copy($ws.Range("A125:V125"))
$ws.Cells.Item(125, 1).Value = 124
$ws.Cells.Item(125, 5).Value = 45263.54166666666
$ws.Cells.Item(125, 6).Value = "Cukaricki"
$ws.Cells.Item(125, 7).Value = 4
$ws.Cells.Item(125, 8).Value = "Radnicki 1923"
$ws.Cells.Item(125, 9).Value = 1
$ws.Cells.Item(125, 10).Value = 1.54
$ws.Cells.Item(125, 11).Value = "01/12/2023 18:43"
$ws.Cells.Item(125, 12).Value = 1.78
$ws.Cells.Item(125, 13).Value = "03/12/2023 12:55"
$ws.Cells.Item(125, 14).Value = 3.76
$ws.Cells.Item(125, 15).Value = "01/12/2023 18:43"
$ws.Cells.Item(125, 16).Value = 3.83
$ws.Cells.Item(125, 17).Value = "03/12/2023 12:55"
$ws.Cells.Item(125, 18).Value = 5.16
$ws.Cells.Item(125, 19).Value = "01/12/2023 18:43"
$ws.Cells.Item(125, 20).Value = 4.01
$ws.Cells.Item(125, 21).Value = "03/12/2023 12:55"
$ws.Cells.Item(125, 22).Value = "https://www.betexplorer.com/football/serbia/super-liga/cukaricki-radnicki-1923/EskpT1XH/"

$ws.Range("A124:V124").Copy($ws.Range("A126:V126"))
$ws.Cells.Item(126, 1).Value = 125
$ws.Cells.Item(126, 5).Value = 45263.625
$ws.Cells.Item(126, 6).Value = "Napredak"
$ws.Cells.Item(126, 7).Value = 0
$ws.Cells.Item(126, 8).Value = "Crvena zvezda"
$ws.Cells.Item(126, 9).Value = 1
$ws.Cells.Item(126, 10).Value = 11.4
$ws.Cells.Item(126, 11).Value = "01/12/2023 18:43"
$ws.Cells.Item(126, 12).Value = 12.86
$ws.Cells.Item(126, 13).Value = "03/12/2023 14:56"
$ws.Cells.Item(126, 14).Value = 6.23
$ws.Cells.Item(126, 15).Value = "01/12/2023 18:43"
$ws.Cells.Item(126, 16).Value = 6.34
$ws.Cells.Item(126, 17).Value = "03/12/2023 14:56"
$ws.Cells.Item(126, 18).Value = 1.17
$ws.Cells.Item(126, 19).Value = "01/12/2023 18:43"
$ws.Cells.Item(126, 20).Value = 1.19
$ws.Cells.Item(126, 21).Value = "03/12/2023 14:56"
$ws.Cells.Item(126, 22).Value = "https://www.betexplorer.com/football/serbia/super-liga/napredak-crvena-zvezda/IDl9XGO3/"

$ws.Range("A124:V124").Copy($ws.Range("A127:V127"))
$ws.Cells.Item(127, 1).Value = 126
$ws.Cells.Item(127, 5).Value = 45263.70833333334
$ws.Cells.Item(127, 6).Value = "TSC"
$ws.Cells.Item(127, 7).Value = 1
$ws.Cells.Item(127, 8).Value = "Radnicki Nis"
$ws.Cells.Item(127, 9).Value = 0
$ws.Cells.Item(127, 10).Value = 1.49
$ws.Cells.Item(127, 11).Value = "01/12/2023 18:43"
$ws.Cells.Item(127, 12).Value = 1.47
$ws.Cells.Item(127, 13).Value = "03/12/2023 16:55"
$ws.Cells.Item(127, 14).Value = 3.88
$ws.Cells.Item(127, 15).Value = "01/12/2023 18:43"
$ws.Cells.Item(127, 16).Value = 3.78
$ws.Cells.Item(127, 17).Value = "03/12/2023 16:55"
$ws.Cells.Item(127, 18).Value = 5.68
$ws.Cells.Item(127, 19).Value = "01/12/2023 18:43"
$ws.Cells.Item(127, 20).Value = 8.09
$ws.Cells.Item(127, 21).Value = "03/12/2023 16:55"
$ws.Cells.Item(127, 22).Value = "https://www.betexplorer.com/football/serbia/super-liga/tsc-backa-topola-radnicki-nis/zijtUsIB/"

$ws.Range("A124:V124").Copy($ws.Range("A128:V128"))
$ws.Cells.Item(128, 1).Value = 127
$ws.Cells.Item(128, 5).Value = 45269.54166666666
$ws.Cells.Item(128, 6).Value = "Radnicki 1923"
$ws.Cells.Item(128, 7).Value = 0
$ws.Cells.Item(128, 8).Value = "TSC"
$ws.Cells.Item(128, 9).Value = 0
$ws.Cells.Item(128, 10).Value = 3.89
$ws.Cells.Item(128, 11).Value = "08/12/2023 01:12"
$ws.Cells.Item(128, 12).Value = 3.01
$ws.Cells.Item(128, 13).Value = "09/12/2023 12:57"
$ws.Cells.Item(128, 14).Value = 3.39
$ws.Cells.Item(128, 15).Value = "08/12/2023 01:12"
$ws.Cells.Item(128, 16).Value = 3.47
$ws.Cells.Item(128, 17).Value = "09/12/2023 12:57"
$ws.Cells.Item(128, 18).Value = 1.8
$ws.Cells.Item(128, 19).Value = "08/12/2023 01:12"
$ws.Cells.Item(128, 20).Value = 2.22
$ws.Cells.Item(128, 21).Value = "09/12/2023 12:57"
$ws.Cells.Item(128, 22).Value = "https://www.betexplorer.com/football/serbia/super-liga/radnicki-1923-tsc-backa-topola/zJQaNW0d/"

$ws.Range("A124:V124").Copy($ws.Range("A129:V129"))
$ws.Cells.Item(129, 1).Value = 128
$ws.Cells.Item(129, 5).Value = 45269.66666666666
$ws.Cells.Item(129, 6).Value = "Crvena zvezda"
$ws.Cells.Item(129, 7).Value = 3
$ws.Cells.Item(129, 8).Value = "Mladost"
$ws.Cells.Item(129, 9).Value = 1
$ws.Cells.Item(129, 10).Value = 1.04
$ws.Cells.Item(129, 11).Value = "08/12/2023 04:12"
$ws.Cells.Item(129, 12).Value = 1.03
$ws.Cells.Item(129, 13).Value = "09/12/2023 15:54"
$ws.Cells.Item(129, 14).Value = 12.91
$ws.Cells.Item(129, 15).Value = "08/12/2023 04:12"
$ws.Cells.Item(129, 16).Value = 17.39
$ws.Cells.Item(129, 17).Value = "09/12/2023 15:54"
$ws.Cells.Item(129, 18).Value = 19.29
$ws.Cells.Item(129, 19).Value = "08/12/2023 04:12"
$ws.Cells.Item(129, 20).Value = 35.04
$ws.Cells.Item(129, 21).Value = "09/12/2023 15:54"
$ws.Cells.Item(129, 22).Value = "https://www.betexplorer.com/football/serbia/super-liga/crvena-zvezda-mladost-lucani/4UnHVfgG/"

$ws.Range("A124:V124").Copy($ws.Range("A130:V130"))
$ws.Cells.Item(130, 1).Value = 129
$ws.Cells.Item(130, 5).Value = 45269.77083333334
$ws.Cells.Item(130, 6).Value = "Zeleznicar Pancevo"
$ws.Cells.Item(130, 7).Value = 2
$ws.Cells.Item(130, 8).Value = "Sp. Subotica"
$ws.Cells.Item(130, 9).Value = 1
$ws.Cells.Item(130, 10).Value = 2.09
$ws.Cells.Item(130, 11).Value = "08/12/2023 06:42"
$ws.Cells.Item(130, 12).Value = 2.35
$ws.Cells.Item(130, 13).Value = "09/12/2023 18:29"
$ws.Cells.Item(130, 14).Value = 3.12
$ws.Cells.Item(130, 15).Value = "08/12/2023 06:42"
$ws.Cells.Item(130, 16).Value = 3.38
$ws.Cells.Item(130, 17).Value = "09/12/2023 18:29"
$ws.Cells.Item(130, 18).Value = 3.25
$ws.Cells.Item(130, 19).Value = "08/12/2023 06:42"
$ws.Cells.Item(130, 20).Value = 2.86
$ws.Cells.Item(130, 21).Value = "09/12/2023 18:29"
$ws.Cells.Item(130, 22).Value = "https://www.betexplorer.com/football/serbia/super-liga/zeleznicar-pancevo-spartak-subotica/rLoLUE8M/"

$ws.Range("A124:V124").Copy($ws.Range("A131:V131"))
$ws.Cells.Item(131, 1).Value = 130
$ws.Cells.Item(131, 5).Value = 45270.54166666666
$ws.Cells.Item(131, 6).Value = "Radnik"
$ws.Cells.Item(131, 7).Value = 1
$ws.Cells.Item(131, 8).Value = "Vozdovac"
$ws.Cells.Item(131, 9).Value = 1
$ws.Cells.Item(131, 10).Value = 2.43
$ws.Cells.Item(131, 11).Value = "08/12/2023 07:12"
$ws.Cells.Item(131, 12).Value = 2.5
$ws.Cells.Item(131, 13).Value = "10/12/2023 12:55"
$ws.Cells.Item(131, 14).Value = 2.92
$ws.Cells.Item(131, 15).Value = "08/12/2023 07:12"
$ws.Cells.Item(131, 16).Value = 2.98
$ws.Cells.Item(131, 17).Value = "10/12/2023 12:55"
$ws.Cells.Item(131, 18).Value = 2.84
$ws.Cells.Item(131, 19).Value = "08/12/2023 07:12"
$ws.Cells.Item(131, 20).Value = 2.98
$ws.Cells.Item(131, 21).Value = "10/12/2023 12:53"
$ws.Cells.Item(131, 22).Value = "https://www.betexplorer.com/football/serbia/super-liga/radnik-surdulica-fk-vozdovac/QqyQTYNS/"

$ws.Range("A124:V124").Copy($ws.Range("A132:V132"))
$ws.Cells.Item(132, 1).Value = 131
$ws.Cells.Item(132, 5).Value = 45270.625
$ws.Cells.Item(132, 6).Value = "Javor"
$ws.Cells.Item(132, 7).Value = 1
$ws.Cells.Item(132, 8).Value = "Cukaricki"
$ws.Cells.Item(132, 9).Value = 2
$ws.Cells.Item(132, 10).Value = 3.45
$ws.Cells.Item(132, 11).Value = "08/12/2023 07:12"
$ws.Cells.Item(132, 12).Value = 3.96
$ws.Cells.Item(132, 13).Value = "10/12/2023 14:30"
$ws.Cells.Item(132, 14).Value = 3.2
$ws.Cells.Item(132, 15).Value = "08/12/2023 07:12"
$ws.Cells.Item(132, 16).Value = 3.36
$ws.Cells.Item(132, 17).Value = "10/12/2023 14:30"
$ws.Cells.Item(132, 18).Value = 1.99
$ws.Cells.Item(132, 19).Value = "08/12/2023 07:12"
$ws.Cells.Item(132, 20).Value = 1.92
$ws.Cells.Item(132, 21).Value = "10/12/2023 14:30"
$ws.Cells.Item(132, 22).Value = "https://www.betexplorer.com/football/serbia/super-liga/javor-cukaricki/dSReOCoj/"

$ws.Range("A124:V124").Copy($ws.Range("A133:V133"))
$ws.Cells.Item(133, 1).Value = 132
$ws.Cells.Item(133, 5).Value = 45270.70833333334
$ws.Cells.Item(133, 6).Value = "IMT Novi Beograd"
$ws.Cells.Item(133, 7).Value = 1
$ws.Cells.Item(133, 8).Value = "Novi Pazar"
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 10).Value = 2.39
$ws.Cells.Item(133, 11).Value = "08/12/2023 07:12"
$ws.Cells.Item(133, 12).Value = 2.6
$ws.Cells.Item(133, 13).Value = "10/12/2023 16:57"
$ws.Cells.Item(133, 14).Value = 3.06
$ws.Cells.Item(133, 15).Value = "08/12/2023 07:12"
$ws.Cells.Item(133, 16).Value = 3.25
$ws.Cells.Item(133, 17).Value = "10/12/2023 16:50"
$ws.Cells.Item(133, 18).Value = 2.78
$ws.Cells.Item(133, 19).Value = "08/12/2023 07:12"
$ws.Cells.Item(133, 20).Value = 2.63
$ws.Cells.Item(133, 21).Value = "10/12/2023 16:57"
$ws.Cells.Item(133, 22).Value = "https://www.betexplorer.com/football/serbia/super-liga/imt-novi-beograd-novi-pazar/25SiPhWq/"

$ws.Range("A124:V124").Copy($ws.Range("A134:V134"))
$ws.Cells.Item(134, 1).Value = 133
$ws.Cells.Item(134, 5).Value = 45271.66666666666
$ws.Cells.Item(134, 6).Value = "Radnicki Nis"
$ws.Cells.Item(134, 7).Value = 2
$ws.Cells.Item(134, 8).Value = "Partizan"
$ws.Cells.Item(134, 9).Value = 1
$ws.Cells.Item(134, 10).Value = 4.95
$ws.Cells.Item(134, 11).Value = "08/12/2023 07:12"
$ws.Cells.Item(134, 12).Value = 5.43
$ws.Cells.Item(134, 13).Value = "11/12/2023 15:58"
$ws.Cells.Item(134, 14).Value = 3.6
$ws.Cells.Item(134, 15).Value = "08/12/2023 07:12"
$ws.Cells.Item(134, 16).Value = 4.24
$ws.Cells.Item(134, 17).Value = "11/12/2023 15:58"
$ws.Cells.Item(134, 18).Value = 1.6
$ws.Cells.Item(134, 19).Value = "08/12/2023 07:12"
$ws.Cells.Item(134, 20).Value = 1.54
$ws.Cells.Item(134, 21).Value = "11/12/2023 15:58"
$ws.Cells.Item(134, 22).Value = "https://www.betexplorer.com/football/serbia/super-liga/radnicki-nis-partizan/CzArnGwc/"

$ws.Range("A124:V124").Copy($ws.Range("A135:V135"))
$ws.Cells.Item(135, 1).Value = 134
$ws.Cells.Item(135, 5).Value = 45271.77083333334
$ws.Cells.Item(135, 6).Value = "Vojvodina"
$ws.Cells.Item(135, 7).Value = 2
$ws.Cells.Item(135, 8).Value = "Napredak"
$ws.Cells.Item(135, 9).Value = 0
$ws.Cells.Item(135, 10).Value = 1.51
$ws.Cells.Item(135, 11).Value = "08/12/2023 07:12"
$ws.Cells.Item(135, 12).Value = 1.47
$ws.Cells.Item(135, 13).Value = "11/12/2023 18:29"
$ws.Cells.Item(135, 14).Value = 3.85
$ws.Cells.Item(135, 15).Value = "08/12/2023 07:12"
$ws.Cells.Item(135, 16).Value = 4.28
$ws.Cells.Item(135, 17).Value = "11/12/2023 18:29"
$ws.Cells.Item(135, 18).Value = 5.46
$ws.Cells.Item(135, 19).Value = "08/12/2023 07:12"
$ws.Cells.Item(135, 20).Value = 6.42
$ws.Cells.Item(135, 21).Value = "11/12/2023 18:29"
$ws.Cells.Item(135, 22).Value = "https://www.betexplorer.com/football/serbia/super-liga/vojvodina-napredak/d4mDWzvA/"

$ws.Range("A124:V124").Copy($ws.Range("A136:V136"))
$ws.Cells.Item(136, 1).Value = 135
$ws.Cells.Item(136, 5).Value = 45275.66666666666
$ws.Cells.Item(136, 6).Value = "Novi Pazar"
$ws.Cells.Item(136, 7).Value = 2
$ws.Cells.Item(136, 8).Value = "Radnik"
$ws.Cells.Item(136, 9).Value = 1
$ws.Cells.Item(136, 10).Value = 1.79
$ws.Cells.Item(136, 11).Value = "13/12/2023 16:12"
$ws.Cells.Item(136, 12).Value = 1.63
$ws.Cells.Item(136, 13).Value = "15/12/2023 15:57"
$ws.Cells.Item(136, 14).Value = 3.22
$ws.Cells.Item(136, 15).Value = "13/12/2023 16:12"
$ws.Cells.Item(136, 16).Value = 3.18
$ws.Cells.Item(136, 17).Value = "15/12/2023 15:57"
$ws.Cells.Item(136, 18).Value = 4.2
$ws.Cells.Item(136, 19).Value = "13/12/2023 16:12"
$ws.Cells.Item(136, 20).Value = 7.05
$ws.Cells.Item(136, 21).Value = "15/12/2023 15:57"
$ws.Cells.Item(136, 22).Value = "https://www.betexplorer.com/football/serbia/super-liga/novi-pazar-radnik-surdulica/4lF2shhS/"

$ws.Range("A124:V124").Copy($ws.Range("A137:V137"))
$ws.Cells.Item(137, 1).Value = 136
$ws.Cells.Item(137, 5).Value = 45275.75
$ws.Cells.Item(137, 6).Value = "Vozdovac"
$ws.Cells.Item(137, 7).Value = 5
$ws.Cells.Item(137, 8).Value = "Zeleznicar Pancevo"
$ws.Cells.Item(137, 9).Value = 1
$ws.Cells.Item(137, 10).Value = 2.17
$ws.Cells.Item(137, 11).Value = "13/12/2023 18:12"
$ws.Cells.Item(137, 12).Value = 1.93
$ws.Cells.Item(137, 13).Value = "15/12/2023 17:58"
$ws.Cells.Item(137, 14).Value = 3.07
$ws.Cells.Item(137, 15).Value = "13/12/2023 18:12"
$ws.Cells.Item(137, 16).Value = 3.48
$ws.Cells.Item(137, 17).Value = "15/12/2023 17:58"
$ws.Cells.Item(137, 18).Value = 3.13
$ws.Cells.Item(137, 19).Value = "13/12/2023 18:12"
$ws.Cells.Item(137, 20).Value = 3.74
$ws.Cells.Item(137, 21).Value = "15/12/2023 17:58"
$ws.Cells.Item(137, 22).Value = "https://www.betexplorer.com/football/serbia/super-liga/fk-vozdovac-zeleznicar-pancevo/dt7UxC0q/"

$ws.Range("A124:V124").Copy($ws.Range("A138:V138"))
$ws.Cells.Item(138, 1).Value = 137
$ws.Cells.Item(138, 5).Value = 45276.58333333334
$ws.Cells.Item(138, 6).Value = "Mladost"
$ws.Cells.Item(138, 7).Value = 1
$ws.Cells.Item(138, 8).Value = "Napredak"
$ws.Cells.Item(138, 9).Value = 2
$ws.Cells.Item(138, 10).Value = 2.16
$ws.Cells.Item(138, 11).Value = "14/12/2023 09:13"
$ws.Cells.Item(138, 12).Value = 2.31
$ws.Cells.Item(138, 13).Value = "16/12/2023 13:58"
$ws.Cells.Item(138, 14).Value = 2.99
$ws.Cells.Item(138, 15).Value = "14/12/2023 09:13"
$ws.Cells.Item(138, 16).Value = 3.08
$ws.Cells.Item(138, 17).Value = "16/12/2023 13:58"
$ws.Cells.Item(138, 18).Value = 3.24
$ws.Cells.Item(138, 19).Value = "14/12/2023 09:13"
$ws.Cells.Item(138, 20).Value = 3.18
$ws.Cells.Item(138, 21).Value = "16/12/2023 13:58"
$ws.Cells.Item(138, 22).Value = "https://www.betexplorer.com/football/serbia/super-liga/mladost-lucani-napredak/IeewyjVe/"

$ws.Range("A124:V124").Copy($ws.Range("A139:V139"))
$ws.Cells.Item(139, 1).Value = 138
$ws.Cells.Item(139, 5).Value = 45276.58333333334
$ws.Cells.Item(139, 6).Value = "Radnicki Nis"
$ws.Cells.Item(139, 7).Value = 0
$ws.Cells.Item(139, 8).Value = "Vojvodina"
$ws.Cells.Item(139, 9).Value = 1
$ws.Cells.Item(139, 10).Value = 2.83
$ws.Cells.Item(139, 11).Value = "14/12/2023 09:13"
$ws.Cells.Item(139, 12).Value = 2.73
$ws.Cells.Item(139, 13).Value = "16/12/2023 13:58"
$ws.Cells.Item(139, 14).Value = 3.11
$ws.Cells.Item(139, 15).Value = "14/12/2023 09:13"
$ws.Cells.Item(139, 16).Value = 3.37
$ws.Cells.Item(139, 17).Value = "16/12/2023 13:57"
$ws.Cells.Item(139, 18).Value = 2.32
$ws.Cells.Item(139, 19).Value = "14/12/2023 09:13"
$ws.Cells.Item(139, 20).Value = 2.45
$ws.Cells.Item(139, 21).Value = "16/12/2023 13:58"
$ws.Cells.Item(139, 22).Value = "https://www.betexplorer.com/football/serbia/super-liga/radnicki-nis-vojvodina/jR8nozg3/"

$ws.Range("A124:V124").Copy($ws.Range("A140:V140"))
$ws.Cells.Item(140, 1).Value = 139
$ws.Cells.Item(140, 5).Value = 45276.66666666666
$ws.Cells.Item(140, 6).Value = "Sp. Subotica"
$ws.Cells.Item(140, 7).Value = 1
$ws.Cells.Item(140, 8).Value = "Crvena zvezda"
$ws.Cells.Item(140, 9).Value = 4
$ws.Cells.Item(140, 10).Value = 11.03
$ws.Cells.Item(140, 11).Value = "14/12/2023 09:13"
$ws.Cells.Item(140, 12).Value = 17.65
$ws.Cells.Item(140, 13).Value = "16/12/2023 15:34"
$ws.Cells.Item(140, 14).Value = 6.27
$ws.Cells.Item(140, 15).Value = "14/12/2023 09:13"
$ws.Cells.Item(140, 16).Value = 7.8
$ws.Cells.Item(140, 17).Value = "16/12/2023 15:34"
$ws.Cells.Item(140, 18).Value = 1.17
$ws.Cells.Item(140, 19).Value = "14/12/2023 09:13"
$ws.Cells.Item(140, 20).Value = 1.13
$ws.Cells.Item(140, 21).Value = "16/12/2023 15:33"
$ws.Cells.Item(140, 22).Value = "https://www.betexplorer.com/football/serbia/super-liga/spartak-subotica-crvena-zvezda/zX5YyWFk/"

$ws.Range("A124:V124").Copy($ws.Range("A141:V141"))
$ws.Cells.Item(141, 1).Value = 140
$ws.Cells.Item(141, 5).Value = 45276.77083333334
$ws.Cells.Item(141, 6).Value = "Partizan"
$ws.Cells.Item(141, 7).Value = 3
$ws.Cells.Item(141, 8).Value = "Radnicki 1923"
$ws.Cells.Item(141, 9).Value = 3
$ws.Cells.Item(141, 10).Value = 1.23
$ws.Cells.Item(141, 11).Value = "14/12/2023 09:13"
$ws.Cells.Item(141, 12).Value = 1.32
$ws.Cells.Item(141, 13).Value = "16/12/2023 18:21"
$ws.Cells.Item(141, 14).Value = 5.38
$ws.Cells.Item(141, 15).Value = "14/12/2023 09:13"
$ws.Cells.Item(141, 16).Value = 4.83
$ws.Cells.Item(141, 17).Value = "16/12/2023 18:29"
$ws.Cells.Item(141, 18).Value = 8.92
$ws.Cells.Item(141, 19).Value = "14/12/2023 09:13"
$ws.Cells.Item(141, 20).Value = 9.67
$ws.Cells.Item(141, 21).Value = "16/12/2023 18:29"
$ws.Cells.Item(141, 22).Value = "https://www.betexplorer.com/football/serbia/super-liga/partizan-radnicki-1923/pl8jpf89/"

$ws.Range("A124:V124").Copy($ws.Range("A142:V142"))
$ws.Cells.Item(142, 1).Value = 141
$ws.Cells.Item(142, 5).Value = 45278.66666666666
$ws.Cells.Item(142, 6).Value = "TSC"
$ws.Cells.Item(142, 7).Value = 3
$ws.Cells.Item(142, 8).Value = "Javor"
$ws.Cells.Item(142, 9).Value = 0
$ws.Cells.Item(142, 10).Value = 1.4
$ws.Cells.Item(142, 11).Value = "14/12/2023 09:13"
$ws.Cells.Item(142, 12).Value = 1.28
$ws.Cells.Item(142, 13).Value = "18/12/2023 15:59"
$ws.Cells.Item(142, 14).Value = 4.22
$ws.Cells.Item(142, 15).Value = "14/12/2023 09:13"
$ws.Cells.Item(142, 16).Value = 5.22
$ws.Cells.Item(142, 17).Value = "18/12/2023 15:59"
$ws.Cells.Item(142, 18).Value = 6.51
$ws.Cells.Item(142, 19).Value = "14/12/2023 09:13"
$ws.Cells.Item(142, 20).Value = 10.42
$ws.Cells.Item(142, 21).Value = "18/12/2023 15:59"
$ws.Cells.Item(142, 22).Value = "https://www.betexplorer.com/football/serbia/super-liga/tsc-backa-topola-javor/rZQgqENF/"

$ws.Range("A124:V124").Copy($ws.Range("A143:V143"))
$ws.Cells.Item(143, 1).Value = 142
$ws.Cells.Item(143, 5).Value = 45278.75
$ws.Cells.Item(143, 6).Value = "Cukaricki"
$ws.Cells.Item(143, 7).Value = 0
$ws.Cells.Item(143, 8).Value = "IMT Novi Beograd"
$ws.Cells.Item(143, 9).Value = 1
$ws.Cells.Item(143, 10).Value = 1.49
$ws.Cells.Item(143, 11).Value = "14/12/2023 09:13"
$ws.Cells.Item(143, 12).Value = 1.51
$ws.Cells.Item(143, 13).Value = "18/12/2023 17:59"
$ws.Cells.Item(143, 14).Value = 3.96
$ws.Cells.Item(143, 15).Value = "14/12/2023 09:13"
$ws.Cells.Item(143, 16).Value = 4.11
$ws.Cells.Item(143, 17).Value = "18/12/2023 17:59"
$ws.Cells.Item(143, 18).Value = 5.56
$ws.Cells.Item(143, 19).Value = "14/12/2023 09:13"
$ws.Cells.Item(143, 20).Value = 6.06
$ws.Cells.Item(143, 21).Value = "18/12/2023 17:59"
$ws.Cells.Item(143, 22).Value = "https://www.betexplorer.com/football/serbia/super-liga/cukaricki-imt-novi-beograd/QcGbrYxM/"
